$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.444287666666667
$ws.Range("H2").Value = 4.332863
$ws.Range("I2").Value = 0.006189216566550864
$ws.Range("J2").Value = 0.006202528009901729
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.506715
$ws.Range("N2").Value = 1.520145
$ws.Range("O2").Value = 0.003122343715987576
$ws.Range("P2").Value = 0.003132472094339857
$ws.Range("Q2").Value = 0.7318422250150001
$ws.Range("R2").Value = 6.586580025135
$ws.Range("S2").Value = [double]"1.932486145345629E-05"
$ws.Range("T2").Value = [double]"1.942924590537849E-05"

$ws.Range("G3").Value = 1.444287666666667
$ws.Range("H3").Value = 4.332863
$ws.Range("I3").Value = 0.006189216566550864
$ws.Range("J3").Value = 0.006202528009901729
$ws.Range("M3").Value = 88.13219433333332
$ws.Range("N3").Value = 264.396583
$ws.Range("O3").Value = 0.5430646480820168
$ws.Range("P3").Value = 0.5448262620252092
$ws.Range("Q3").Value = 127.2882413119032
$ws.Range("R3").Value = 1145.594171807129
$ws.Range("S3").Value = 0.003361144716617333
$ws.Range("T3").Value = 0.003379300150741419

$ws.Range("G4").Value = 1.444287666666667
$ws.Range("H4").Value = 4.332863
$ws.Range("I4").Value = 0.006189216566550864
$ws.Range("J4").Value = 0.006202528009901729
$ws.Range("M4").Value = 1.5741895
$ws.Range("N4").Value = 3.148379
$ws.Range("O4").Value = 0.009700049718478087
$ws.Range("P4").Value = 0.006487676741301404
$ws.Range("Q4").Value = 2.273582479846167
$ws.Range("R4").Value = 13.641494879077
$ws.Range("S4").Value = [double]"6.003570841397161E-05"
$ws.Range("T4").Value = [double]"4.023999670710994E-05"

$ws.Range("G5").Value = 1.444287666666667
$ws.Range("H5").Value = 4.332863
$ws.Range("I5").Value = 0.006189216566550864
$ws.Range("J5").Value = 0.006202528009901729
$ws.Range("M5").Value = 72.07364666666666
$ws.Range("N5").Value = 216.22094
$ws.Range("O5").Value = 0.4441129584835175
$ws.Range("P5").Value = 0.4455535891391496
$ws.Range("Q5").Value = 104.0950789723578
$ws.Range("R5").Value = 936.8557107512198
$ws.Range("S5").Value = 0.002748711280066102
$ws.Range("T5").Value = 0.002763558616547822

$ws.Range("I6").Value = 0.003021900187532335
$ws.Range("J6").Value = 0.003028399532437424
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.506715
$ws.Range("N6").Value = 1.520145
$ws.Range("O6").Value = 0.003122343715987576
$ws.Range("P6").Value = 0.003132472094339857
$ws.Range("Q6").Value = 0.357323763555
$ws.Range("R6").Value = 3.215913871995
$ws.Range("S6").Value = [double]"9.435411060883264E-06"
$ws.Range("T6").Value = [double]"9.486377025872101E-06"

$ws.Range("I7").Value = 0.003021900187532335
$ws.Range("J7").Value = 0.003028399532437424
$ws.Range("M7").Value = 88.13219433333332
$ws.Range("N7").Value = 264.396583
$ws.Range("O7").Value = 0.5430646480820168
$ws.Range("P7").Value = 0.5448262620252092
$ws.Range("Q7").Value = 62.14879640339699
$ws.Range("R7").Value = 559.3391676305729
$ws.Range("S7").Value = 0.001641087161881228
$ws.Range("T7").Value = 0.001649951597176773

$ws.Range("I8").Value = 0.003021900187532335
$ws.Range("J8").Value = 0.003028399532437424
$ws.Range("M8").Value = 1.5741895
$ws.Range("N8").Value = 3.148379
$ws.Range("O8").Value = 0.009700049718478087
$ws.Range("P8").Value = 0.006487676741301404
$ws.Range("Q8").Value = 1.1100822290415
$ws.Range("R8").Value = 6.660493374248999
$ws.Range("S8").Value = [double]"2.93125820633419E-05"
$ws.Range("T8").Value = [double]"1.964727720996232E-05"

$ws.Range("I9").Value = 0.003021900187532335
$ws.Range("J9").Value = 0.003028399532437424
$ws.Range("M9").Value = 72.07364666666666
$ws.Range("N9").Value = 216.22094
$ws.Range("O9").Value = 0.4441129584835175
$ws.Range("P9").Value = 0.4455535891391496
$ws.Range("Q9").Value = 50.82467793546
$ws.Range("R9").Value = 457.4221014191399
$ws.Range("S9").Value = 0.001342065032526882
$ws.Range("T9").Value = 0.001349314281024817

$ws.Range("G10").Value = 123.254125
$ws.Range("H10").Value = 369.762375
$ws.Range("I10").Value = 0.528181808895687
$ws.Range("J10").Value = 0.5293177947110922
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.506715
$ws.Range("N10").Value = 1.520145
$ws.Range("O10").Value = 0.003122343715987576
$ws.Range("P10").Value = 0.003132472094339857
$ws.Range("Q10").Value = 62.45471394937501
$ws.Range("R10").Value = 562.092425544375
$ws.Range("S10").Value = 0.001649165151904399
$ws.Range("T10").Value = 0.001658073220970009

$ws.Range("G11").Value = 123.254125
$ws.Range("H11").Value = 369.762375
$ws.Range("I11").Value = 0.528181808895687
$ws.Range("J11").Value = 0.5293177947110922
$ws.Range("M11").Value = 88.13219433333332
$ws.Range("N11").Value = 264.396583
$ws.Range("O11").Value = 0.5430646480820168
$ws.Range("P11").Value = 0.5448262620252092
$ws.Range("Q11").Value = 10862.65649688496
$ws.Range("R11").Value = 97763.90847196462
$ws.Range("S11").Value = 0.2868368681712593
$ws.Range("T11").Value = 0.2883862355158714

$ws.Range("G12").Value = 123.254125
$ws.Range("H12").Value = 369.762375
$ws.Range("I12").Value = 0.528181808895687
$ws.Range("J12").Value = 0.5293177947110922
$ws.Range("M12").Value = 1.5741895
$ws.Range("N12").Value = 3.148379
$ws.Range("O12").Value = 0.009700049718478087
$ws.Range("P12").Value = 0.006487676741301404
$ws.Range("Q12").Value = 194.0253494066875
$ws.Range("R12").Value = 1164.152096440125
$ws.Range("S12").Value = 0.005123389806683855
$ws.Range("T12").Value = 0.003434042745504104

$ws.Range("G13").Value = 123.254125
$ws.Range("H13").Value = 369.762375
$ws.Range("I13").Value = 0.528181808895687
$ws.Range("J13").Value = 0.5293177947110922
$ws.Range("M13").Value = 72.07364666666666
$ws.Range("N13").Value = 216.22094
$ws.Range("O13").Value = 0.4441129584835175
$ws.Range("P13").Value = 0.4455535891391496
$ws.Range("Q13").Value = 8883.374255459166
$ws.Range("R13").Value = 79950.3682991325
$ws.Range("S13").Value = 0.2345723857658394
$ws.Range("T13").Value = 0.2358394432287467

$ws.Range("G14").Value = 1.502435
$ws.Range("H14").Value = 3.00487
$ws.Range("I14").Value = 0.006438395761993292
$ws.Range("J14").Value = 0.00430149541795192
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 0.506715
$ws.Range("N14").Value = 1.520145
$ws.Range("O14").Value = 0.003122343715987576
$ws.Range("P14").Value = 0.003132472094339857
$ws.Range("Q14").Value = 0.761306351025
$ws.Range("R14").Value = 4.56783810615
$ws.Range("S14").Value = [double]"2.01028845485008E-05"
$ws.Range("T14").Value = [double]"1.347431436066515E-05"

$ws.Range("G15").Value = 1.502435
$ws.Range("H15").Value = 3.00487
$ws.Range("I15").Value = 0.006438395761993292
$ws.Range("J15").Value = 0.00430149541795192
$ws.Range("M15").Value = 88.13219433333332
$ws.Range("N15").Value = 264.396583
$ws.Range("O15").Value = 0.5430646480820168
$ws.Range("P15").Value = 0.5448262620252092
$ws.Range("Q15").Value = 132.4128933932016
$ws.Range("R15").Value = 794.4773603592099
$ws.Range("S15").Value = 0.003496465128699635
$ws.Range("T15").Value = 0.00234356766968131

$ws.Range("G16").Value = 1.502435
$ws.Range("H16").Value = 3.00487
$ws.Range("I16").Value = 0.006438395761993292
$ws.Range("J16").Value = 0.00430149541795192
$ws.Range("M16").Value = 1.5741895
$ws.Range("N16").Value = 3.148379
$ws.Range("O16").Value = 0.009700049718478087
$ws.Range("P16").Value = 0.006487676741301404
$ws.Range("Q16").Value = 2.3651174014325
$ws.Range("R16").Value = 9.46046960573
$ws.Range("S16").Value = [double]"6.245275899857354E-05"
$ws.Range("T16").Value = [double]"2.790671177586124E-05"

$ws.Range("G17").Value = 1.502435
$ws.Range("H17").Value = 3.00487
$ws.Range("I17").Value = 0.006438395761993292
$ws.Range("J17").Value = 0.00430149541795192
$ws.Range("M17").Value = 72.07364666666666
$ws.Range("N17").Value = 216.22094
$ws.Range("O17").Value = 0.4441129584835175
$ws.Range("P17").Value = 0.4455535891391496
$ws.Range("Q17").Value = 108.2859693296333
$ws.Range("R17").Value = 649.7158159777999
$ws.Range("S17").Value = 0.002859374989746582
$ws.Range("T17").Value = 0.001916546722134084

$ws.Range("G18").Value = 106.4494656666667
$ws.Range("H18").Value = 319.348397
$ws.Range("I18").Value = 0.4561686785882365
$ws.Range("J18").Value = 0.4571497823286167
$ws.Range("K18").Value = 3.0
$ws.Range("L18").Value = 1.0
$ws.Range("M18").Value = 0.506715
$ws.Range("N18").Value = 1.520145
$ws.Range("O18").Value = 0.003122343715987576
$ws.Range("P18").Value = 0.003132472094339857
$ws.Range("Q18").Value = 53.93954099528499
$ws.Range("R18").Value = 485.455868957565
$ws.Range("S18").Value = 0.001424315407020336
$ws.Range("T18").Value = 0.001432008936077932

$ws.Range("G19").Value = 106.4494656666667
$ws.Range("H19").Value = 319.348397
$ws.Range("I19").Value = 0.4561686785882365
$ws.Range("J19").Value = 0.4571497823286167
$ws.Range("M19").Value = 88.13219433333332
$ws.Range("N19").Value = 264.396583
$ws.Range("O19").Value = 0.5430646480820168
$ws.Range("P19").Value = 0.5448262620252092
$ws.Range("Q19").Value = 9381.624994814158
$ws.Range("R19").Value = 84434.62495332744
$ws.Range("S19").Value = 0.2477290829035593
$ws.Range("T19").Value = 0.2490672070917383

$ws.Range("G20").Value = 106.4494656666667
$ws.Range("H20").Value = 319.348397
$ws.Range("I20").Value = 0.4561686785882365
$ws.Range("J20").Value = 0.4571497823286167
$ws.Range("M20").Value = 1.5741895
$ws.Range("N20").Value = 3.148379
$ws.Range("O20").Value = 0.009700049718478087
$ws.Range("P20").Value = 0.006487676741301404
$ws.Range("Q20").Value = 167.5716311330771
$ws.Range("R20").Value = 1005.429786798463
$ws.Range("S20").Value = 0.004424858862318344
$ws.Range("T20").Value = 0.002965840010104367

$ws.Range("G21").Value = 106.4494656666667
$ws.Range("H21").Value = 319.348397
$ws.Range("I21").Value = 0.4561686785882365
$ws.Range("J21").Value = 0.4571497823286167
$ws.Range("M21").Value = 72.07364666666666
$ws.Range("N21").Value = 216.22094
$ws.Range("O21").Value = 0.4441129584835175
$ws.Range("P21").Value = 0.4455535891391496
$ws.Range("Q21").Value = 7672.201176314797
$ws.Range("R21").Value = 69049.81058683318
$ws.Range("S21").Value = 0.2025904214153385
$ws.Range("T21").Value = 0.2036847262906961
